# 05/01/26 - adding Sainsbury's shopping and 200GB SIM card purchase
#
# Updates the January row with the extra Sainsbury's grocery shop and the
# larger "Special" spend for the 200GB SIM card, tidies up the
# "Total + Utilities" column header, appends a note about the SIM card to
# the existing threaded comment on O2, and leaves the selection where the
# author left off (J7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Groceries (I2) and Special (O2) for January go up after the extra
# Sainsbury's shop + SIM card purchase. L2 (=H2+I2+J2+K2) and
# N2 (=B2+L2+M2) are formulas and recalculate automatically.
$ws.Range("I2").Value = 74
$ws.Range("O2").Value = 340

# Column header N1 no longer includes Rent in its label.
$ws.Range("N1").Value = "Total + Utilities"

# Extend the existing threaded/legacy comment on O2 with a mention of the
# large mobile data SIM card purchase.
$comment = $ws.Range("O2").Comment
$existingText = $comment.Text()
$comment.Text($existingText + " + large mobile data SIM card")

# Restore the author's last selection.
$ws.Range("J7").Select()
